$d = $word.ActiveDocument

function Get-GoBackParagraph($doc) {
    foreach ($p in $doc.Paragraphs) {
        foreach ($bm in $p.Range.Bookmarks) {
            if ($bm.Name -eq "_GoBack") {
                return $p
            }
        }
    }
    return $doc.Paragraphs.Last
}

$pkgHead = '<?xml version="1.0"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>'
$pkgTail = '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

# --- Step 1 -----------------------------------------------------------
# Insert a standalone blank paragraph immediately before the paragraph
# that holds the "_GoBack" bookmark (i.e. right after the lone-space
# paragraph that currently precedes it). A run-less <w:p> fragment, when
# inserted via InsertXML, becomes its own new paragraph rather than
# merging into the target paragraph's runs.
$emptyFrag = '<w:p><w:pPr><w:rPr><w:rFonts w:hint="default"/><w:lang w:eastAsia="zh-CN"/></w:rPr></w:pPr></w:p>'

$target = Get-GoBackParagraph($d)
$insertAt = $target.Range.Start
$rng = $d.Range($insertAt, $insertAt)
$rng.InsertXML($pkgHead + $emptyFrag + $pkgTail)

# --- Step 2 -----------------------------------------------------------
# The paragraph just created is now the immediate predecessor of the
# bookmark paragraph. Insert the "V11:" title paragraph plus the
# HYPERLINK field-code paragraph right at its start; since the field
# paragraph is the LAST paragraph of this fragment its runs merge into
# that (still-empty) target paragraph, which keeps the desired
# "fldChar end" / rStyle formatting on the visible URL text intact and
# leaves the title paragraph as its own standalone paragraph ahead of it.
$target2 = Get-GoBackParagraph($d)
$blankPara = $target2.Previous()
$insertAt2 = $blankPara.Range.Start
$rng2 = $d.Range($insertAt2, $insertAt2)

$titleFrag = '<w:p><w:r><w:rPr><w:rFonts w:hint="default"/><w:lang w:eastAsia="zh-CN"/></w:rPr><w:t xml:space="preserve">V11: </w:t></w:r><w:r><w:t>Hướng dẫn debug trong Visual Studio Code</w:t></w:r></w:p>'
$fieldFrag = '<w:p><w:pPr><w:rPr><w:rFonts w:hint="default"/><w:lang w:eastAsia="zh-CN"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:hint="default"/><w:lang w:eastAsia="zh-CN"/></w:rPr><w:fldChar w:fldCharType="begin"/></w:r><w:r><w:rPr><w:rFonts w:hint="default"/><w:lang w:eastAsia="zh-CN"/></w:rPr><w:instrText xml:space="preserve"> HYPERLINK "https://youtu.be/yyA_E9xxlsk" </w:instrText></w:r><w:r><w:rPr><w:rFonts w:hint="default"/><w:lang w:eastAsia="zh-CN"/></w:rPr><w:fldChar w:fldCharType="separate"/></w:r><w:r><w:rPr><w:rStyle w:val="4"/><w:rFonts w:hint="default"/><w:lang w:eastAsia="zh-CN"/></w:rPr><w:t>https://youtu.be/yyA_E9xxlsk</w:t></w:r><w:r><w:rPr><w:rFonts w:hint="default"/><w:lang w:eastAsia="zh-CN"/></w:rPr><w:fldChar w:fldCharType="end"/></w:r></w:p>'

$rng2.InsertXML($pkgHead + $titleFrag + $fieldFrag + $pkgTail)

Write-Output "Inserted V11 title + hyperlink field + trailing blank paragraph."
